$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Fri Jul 28 08:28:30 UTC 2023
$ws.Range("D2").Value = "29.158.51"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.861.24"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.7078"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'240.96"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("D8").Value = "'0.3087"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "'0.07640"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "'24.69"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "'0.08354"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "1.859.80"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "'5.182"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "'0.7077"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "'91.10"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "29.180.38"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "'5.907"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'242.98"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").Value = "'0.000007801"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "2.115.35"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'13.07"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'7.869"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'0.1585"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'163.28"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'8.932"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'18.46"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "'1.338"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'1.498"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'4.397"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'4.217"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "'0.05142"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").Value = "'0.7964"
$ws.Range("E34").Value = "  +9.96%  "
$ws.Range("D35").Value = "'1.912"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'1.162"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.01841"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "'2.694"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "1.168.23"
$ws.Range("E40").Value = "  -6.15%  "
$ws.Range("D41").Value = "'6.194"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'0.8909"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'72.77"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("D45").Value = "'101.93"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "2.012.67"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'0.5195"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'1.769"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.328"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "'1.002"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4276"
$ws.Range("E51").Value = "  -0.97%  "
